$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits at the very end of the
#    document (after "Fluent in Spanish"). Word relocates this
#    bookmark to mark the most-recently-edited spot, so remove the
#    old one now; it will be re-created inside the new paragraph below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Find the "EDUCATION:" heading paragraph (by absolute paragraph
#    index) and insert a brand new, empty paragraph immediately
#    before it.
# ------------------------------------------------------------------
$eduIndex = -1
$i = 0
foreach ($para in $d.Paragraphs) {
    $i++
    $trimmed = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($trimmed -eq "EDUCATION:") {
        $eduIndex = $i
        break
    }
}

$eduPara = $d.Paragraphs($eduIndex)
$eduPara.Range.InsertParagraphBefore()

# The "EDUCATION:" paragraph shifted down by one; the newly created
# (still empty) paragraph now occupies its former slot.
$newPara = $d.Paragraphs($eduIndex)
$targetRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

# ------------------------------------------------------------------
# 3. Fill the new paragraph with the OBJECTIVE text, matching the
#    exact run/formatting layout, and embed the relocated _GoBack
#    bookmark in the middle of it.
# ------------------------------------------------------------------
$objectiveXml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="Heading1"/>
<w:rPr><w:sz w:val="24"/></w:rPr>
</w:pPr>
<w:r><w:t xml:space="preserve">OBJECTIVE: </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Eager to drive </w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">solutions at </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Motorola</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> on a full-time basis</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$targetRange.InsertXML($objectiveXml)

Write-Host "New paragraph text: $($d.Paragraphs($eduIndex).Range.Text)"
